# The workbook gained one new data row: a new record was inserted as row 50
# (pushing the former rows 50-107 down to 51-108, which keeps all of their
# original values intact), and the newly inserted row 50 was populated with
# a fresh observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 50; this shifts rows 50..107 down to 51..108
# and keeps the rest of the sheet (rows 1..49) untouched.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new record's data.
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = 44638
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = 100112052
$ws.Range("G50").Value = "Albahaca"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 120
$ws.Range("K50").Value = 6000
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = 6000
$ws.Range("N50").Value = "`$/docena de matas"
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 1000
$ws.Range("Q50").Value = 6
$ws.Range("R50").Value = "Hortaliza"
